$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# --- Step 1: preserve the footnote row (currently row 107: empty A107 +
#     "※4/8より..." text in B107) by relocating its formatting/value down to
#     row 109 BEFORE row 107 is overwritten with new data. ---
$ws.Range("A107:B107").Copy()
$ws.Range("A109:B109").PasteSpecial(-4122)
$ws.Range("B109").Value = "※4/8より健康相談窓口と帰国者・接触者相談センターを統合"

# --- Step 2: write the two new daily-data rows (5/11 and 5/12) into what
#     used to be the footnote row (107) plus a brand-new row (108). ---
$ws.Range("A107").Value = 43962
$ws.Range("B107").Value = 455
$ws.Range("C107").Value = 36234
$ws.Range("D107").Value = 113
$ws.Range("E107").Value = 7345

$ws.Range("A108").Value = 43963
$ws.Range("B108").Value = 314
$ws.Range("C108").Value = 36548
$ws.Range("D108").Value = 92
$ws.Range("E108").Value = 7437

# --- Step 3: match the formatting of the preceding data row (106) for the
#     two freshly added rows. ---
$ws.Range("A106:E106").Copy()
$ws.Range("A107:E108").PasteSpecial(-4122)

# --- Step 4: grow the printed area by two rows (108 -> 110). ---
$nm = $wb.Names.Item(1)
$nm.RefersTo = '=相談件数!$A$1:$E$110'

# --- Step 5: move the selection to the new last cell, as in the source file. ---
[void]$ws.Range("E109").Select()
